$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number need to be forced back to
# Text format first, otherwise the COM layer auto-converts "231.34" etc. into
# a numeric value (and would also strip things like trailing zeros).
$ws.Range('D2').Value = '36.381.18'
$ws.Range('E2').Value = '  -0.70%  '
$ws.Range('D3').Value = '2.020.38'
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.34'
$ws.Range('E5').Value = '  -9.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.598'
$ws.Range('E6').Value = '  -2.17%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.67'
$ws.Range('E8').Value = '  -1.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.370'
$ws.Range('E9').Value = '  -1.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '57.01'
$ws.Range('E10').Value = '  +3.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0745'
$ws.Range('E11').Value = '  -2.19%  '
$ws.Range('E12').Value = '  -1.21%  '
$ws.Range('D13').Value = '2.323.12'
$ws.Range('E13').Value = '  +1.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.19'
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.97'
$ws.Range('E15').Value = '  -5.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.755'
$ws.Range('E16').Value = '  -4.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.07'
$ws.Range('E17').Value = '  -1.53%  '
$ws.Range('D18').Value = '2.042.89'
$ws.Range('E18').Value = '  +3.20%  '
$ws.Range('D19').Value = '36.597.08'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '67.51'
$ws.Range('E20').Value = '  -3.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.48'
$ws.Range('E21').Value = '  +8.65%  '
$ws.Range('D22').Value = '0.0₃0792'
$ws.Range('E22').Value = '  -3.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '220.44'
$ws.Range('E23').Value = '  -6.32%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.37'
$ws.Range('E25').Value = '  +0.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.35'
$ws.Range('E26').Value = '  -7.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.31'
$ws.Range('E27').Value = '  -0.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.59'
$ws.Range('E28').Value = '  -2.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.128'
$ws.Range('E29').Value = '  +4.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '18.84'
$ws.Range('E30').Value = '  -2.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.34'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('E32').Value = '  -1.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.35'
$ws.Range('E33').Value = '  -3.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0598'
$ws.Range('E34').Value = '  -4.67%  '
$ws.Range('E35').Value = '  +4.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.23'
$ws.Range('E36').Value = '  -2.60%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.76'
$ws.Range('E38').Value = '  -2.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.26'
$ws.Range('E39').Value = '  -4.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.76'
$ws.Range('E40').Value = '  +5.60%  '
$ws.Range('E41').Value = '  -2.62%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.457.87'
$ws.Range('E42').Value = '  +0.99%  '
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0934'
$ws.Range('E43').Value = '  +2.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '93.20'
$ws.Range('E44').Value = '  +5.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0202'
$ws.Range('E45').Value = '  -2.31%  '
$ws.Range('E46').Value = '  -4.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.44'
$ws.Range('E47').Value = '  -0.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.95'
$ws.Range('E48').Value = '  +33.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.997'
$ws.Range('E49').Value = '  -1.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.88'
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.84'
$ws.Range('E51').Value = '  -0.13%  '

Write-Host "done"
